# Fix sch pin error: PA2 --> PA5 (cell B2 on Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "PA5"

# Update the active selection to match the edited cell
$ws.Range("B2").Select()
